$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

# Column A holds a date-style label ("01-07-2021") that must be stored as
# literal text (matching every other row in the "Serie" column), not be
# auto-converted to a date serial number. Briefly marking the cell as Text
# before assigning the value keeps Excel's "looks like a date" smart
# detection from kicking in; switching the cell back to the "Normal" style
# afterwards keeps its formatting identical to the rest of the column.
$cA = $ws.Cells.Item($row, 1)
$cA.NumberFormat = "@"
$cA.Value = "01-07-2021"
$cA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 233155
$ws.Cells.Item($row, 3).Value = 44254
$ws.Cells.Item($row, 4).Value = 3604
$ws.Cells.Item($row, 5).Value = 27333
$ws.Cells.Item($row, 6).Value = 2196
$ws.Cells.Item($row, 7).Value = 98846
$ws.Cells.Item($row, 8).Value = 56921
$ws.Cells.Item($row, 9).Value = 231266
